# "changes in concise marksheet" — update the Corr/total marks on the
# marksheet for roll number 1401ME59 (sheet "quiz"):
#   - B11 (Marking / Right): 3  -> 5
#   - B12 (Total / Right):  42 -> 70
#   - E12 (Total / Max, "correct/total"): "38/84" -> "70/140"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 70
$ws.Range("E12").Value = "70/140"
